# Updates the cryptos price/volume table on Sheet1 to the latest scraped
# values (GitHub Actions refresh). Columns: B=Coin, C=Link, D=Price,
# E=Volume(1h). Price cells whose new value could be mis-parsed as a
# number by Excel (losing trailing zeros / switching to sci. notation)
# are forced to text format first so the literal string is preserved,
# matching how the sheet already stores every Price/Volume cell as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.234.56"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "3.619.89"
$ws.Range("E3").Value = "  +5.98%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.74"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "652.42"
$ws.Range("E6").Value = "  +5.40%  "
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.404"
$ws.Range("E8").Value = "  +3.36%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.993"
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("D11").Value = "3.619.64"
$ws.Range("E11").Value = "  +6.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.54"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.30"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "4.295.16"
$ws.Range("E15").Value = "  +5.73%  "
$ws.Range("D16").Value = "95.219.28"
$ws.Range("E16").Value = "  +2.25%  "
$ws.Range("E17").Value = "  +3.05%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.622.39"
$ws.Range("E18").Value = "  +6.43%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.88"
$ws.Range("E19").Value = "  -4.11%  "
$ws.Range("E20").Value = "  +11.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.97"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.58"
$ws.Range("E22").Value = "  +5.80%  "
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "506.26"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000196"
$ws.Range("E25").Value = "  +6.97%  "
$ws.Range("E26").Value = "  -1.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.89"
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.68"
$ws.Range("E28").Value = "  +5.58%  "
$ws.Range("D29").Value = "3.816.24"
$ws.Range("E29").Value = "  +6.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.11"
$ws.Range("E30").Value = "  +14.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.29"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.991"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "32.11"
$ws.Range("E35").Value = "  +10.98%  "
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.558"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.14"
$ws.Range("E38").Value = "  +9.09%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "572.13"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.47"
$ws.Range("E40").Value = "  +4.89%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.926"
$ws.Range("E42").Value = "  +3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.149"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "35.49"
$ws.Range("E44").Value = "  +46.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.71"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.69"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.66"
$ws.Range("E47").Value = "  +4.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.21"
$ws.Range("E48").Value = "  +5.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0412"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("E50").Value = "  -4.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.67"
$ws.Range("E51").Value = "  +1.14%  "
